# "Added mobile code countries table"
# Insert a MobileCode column on both the English and Arabic sides of the
# Country/State/City table:
#   Old: Country(en) | State(en) | City(en) | Country(ar) | State(ar) | City(ar)
#   New: Country(en) | MobileCode(en) | State(en) | City(en) | Country(ar) | MobileCode(ar) | State(ar) | City(ar)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column right after Country(en) (old column B, "State(en)") -
# this becomes the new MobileCode(en) column.
$ws.Columns("B").Insert()

# Insert a new column right after Country(ar) (now shifted to column F by the
# previous insert) - this becomes the new MobileCode(ar) column.
$ws.Columns("F").Insert()

# Header row
$ws.Range("B1").Value = "MobileCode(en)"
$ws.Range("F1").Value = "MobileCode(ar)"
$ws.Range("B1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true

# Data row
$ws.Range("B2").Value = 91
$ws.Range("F2").Value = 91

# Re-apply column widths as close as this engine's column-width quantization
# allows to the widths recorded after the edit (auto/manual resize following
# the column insert).
$ws.Columns("A").ColumnWidth = 15.833333333333332
$ws.Columns("B").ColumnWidth = 14.166666666666668
$ws.Columns("C").ColumnWidth = 15.666666666666668
$ws.Columns("D").ColumnWidth = 19.333333333333336
$ws.Columns("E").ColumnWidth = 15.666666666666668
$ws.Columns("F").ColumnWidth = 14.833333333333332
$ws.Columns("G").ColumnWidth = 12.333333333333332
$ws.Columns("H").ColumnWidth = 12.5

# Match the saved selection position.
$ws.Range("G5").Select() | Out-Null
